$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 3 and row 5 and must be swapped:
# A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr3 = "{0}3" -f $col
    $addr5 = "{0}5" -f $col
    $v3 = $ws.Range($addr3).Value2
    $v5 = $ws.Range($addr5).Value2
    $ws.Range($addr3).Value = $v5
    $ws.Range($addr5).Value = $v3
}
